# Generate Report for Archive
#
# 1. The "Ready for handoff" status text is updated to "In Translation"
#    everywhere it is used (Overview sheet columns E/F and the per-locale
#    sheets' Status column).
# 2. The width of those same status columns is narrowed to fit the new,
#    shorter text (from ~17.22 to ~13.41 characters).

$wb = $excel.ActiveWorkbook

# --- Sheet 1: Overview ---------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")

# Update status values in columns E (zh-cn) and F (de-de), rows 2-4.
for ($r = 2; $r -le 4; $r++) {
    foreach ($col in @("E", "F")) {
        $cell = $overview.Range("$col$r")
        if ($cell.Value2 -eq "Ready for handoff") {
            $cell.Value = "In Translation"
        }
    }
}

# Narrow columns E and F to match the shorter text.
# (Target OOXML width is ~13.41 characters; this host's ColumnWidth setter
# only lands on a 1/6-character grid, so 12.5 is the input that rounds to
# the closest representable stored width, 13.3333...)
$overview.Range("E1:F1").ColumnWidth = 12.5

# --- Sheet 2: zh-cn -------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

for ($r = 2; $r -le 4; $r++) {
    $cell = $zhcn.Range("C$r")
    if ($cell.Value2 -eq "Ready for handoff") {
        $cell.Value = "In Translation"
    }
}

$zhcn.Range("C1").ColumnWidth = 12.5

# --- Sheet 3: de-de -------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

for ($r = 2; $r -le 4; $r++) {
    $cell = $dede.Range("C$r")
    if ($cell.Value2 -eq "Ready for handoff") {
        $cell.Value = "In Translation"
    }
}

$dede.Range("C1").ColumnWidth = 12.5
